# Scheduled market-data refresh for the Typhon data-center profit sheets.
# Updates the computed price/profit columns (H:N) for the leves whose
# underlying Universalis market prices changed since the last run.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 9025219
$ws.Range("I80").Value = 310
$ws.Range("J80").Value = 14333988
$ws.Range("K80").Value = 930
$ws.Range("L80").Value = 43001964
$ws.Range("M80").Value = 68
$ws.Range("N80").Value = -43003960
$ws.Range("H83").Value = 9025219
$ws.Range("I83").Value = 310
$ws.Range("J83").Value = 14333988
$ws.Range("K83").Value = 2790
$ws.Range("L83").Value = 129005892
$ws.Range("M83").Value = 2202
$ws.Range("N83").Value = -129015876
$ws.Range("H116").Value = 3846.2856
$ws.Range("I116").Value = 2285
$ws.Range("J116").Value = 4807.077
$ws.Range("K116").Value = 2285
$ws.Range("L116").Value = 4807.077
$ws.Range("M116").Value = 1157
$ws.Range("N116").Value = -11691.077
$ws.Range("H132").Value = 2487.9285
$ws.Range("I132").Value = 2823.147
$ws.Range("K132").Value = 8469.440999999999
$ws.Range("M132").Value = -5939.440999999999
$ws.Range("H135").Value = 21742186
$ws.Range("I135").Value = 844.4706
$ws.Range("K135").Value = 7600.2354
$ws.Range("M135").Value = -5065.2354
$ws.Range("H137").Value = 1319.8096
$ws.Range("I137").Value = 1251.0667
$ws.Range("K137").Value = 3753.2001
$ws.Range("M137").Value = -1203.2001
$ws.Range("H138").Value = 35716884
$ws.Range("J138").Value = 3306.2144
$ws.Range("L138").Value = 9918.643199999999
$ws.Range("N138").Value = -20198.6432
$ws.Range("H141").Value = 3404.9285
$ws.Range("I141").Value = 2969.9092
$ws.Range("K141").Value = 8909.7276
$ws.Range("M141").Value = -3729.7276

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5317.5493
$ws.Range("I32").Value = 4500.492
$ws.Range("K32").Value = 4500.492
$ws.Range("M32").Value = -4213.492
$ws.Range("H45").Value = 2317.1428
$ws.Range("I45").Value = 1714.1428
$ws.Range("J45").Value = 3523.1428
$ws.Range("K45").Value = 1714.1428
$ws.Range("L45").Value = 3523.1428
$ws.Range("M45").Value = -1337.1428
$ws.Range("N45").Value = -4277.1428
$ws.Range("H63").Value = 2054.5
$ws.Range("I63").Value = 2094.9
$ws.Range("J63").Value = 1852.5
$ws.Range("K63").Value = 2094.9
$ws.Range("L63").Value = 1852.5
$ws.Range("M63").Value = -1408.9
$ws.Range("N63").Value = -3224.5
$ws.Range("H66").Value = 2054.5
$ws.Range("I66").Value = 2094.9
$ws.Range("J66").Value = 1852.5
$ws.Range("K66").Value = 10474.5
$ws.Range("L66").Value = 9262.5
$ws.Range("M66").Value = -7042.5
$ws.Range("N66").Value = -16126.5
$ws.Range("H74").Value = 52633628
$ws.Range("J74").Value = 2532.0715
$ws.Range("L74").Value = 2532.0715
$ws.Range("N74").Value = -4280.0715
$ws.Range("H77").Value = 52633628
$ws.Range("J77").Value = 2532.0715
$ws.Range("L77").Value = 12660.3575
$ws.Range("N77").Value = -21396.3575
$ws.Range("H122").Value = 1570.1578
$ws.Range("I122").Value = 1420.5667
$ws.Range("J122").Value = 2131.125
$ws.Range("K122").Value = 4261.7001
$ws.Range("L122").Value = 6393.375
$ws.Range("M122").Value = -1811.7001
$ws.Range("N122").Value = -11293.375

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3080.1667
$ws.Range("I31").Value = 2348.1667
$ws.Range("K31").Value = 2348.1667
$ws.Range("M31").Value = -2053.1667
$ws.Range("H34").Value = 3080.1667
$ws.Range("I34").Value = 2348.1667
$ws.Range("K34").Value = 2348.1667
$ws.Range("M34").Value = -2146.1667
$ws.Range("H52").Value = 35666.668
$ws.Range("J52").Value = 35666.668
$ws.Range("L52").Value = 35666.668
$ws.Range("N52").Value = -36254.668
$ws.Range("H86").Value = 29971.6
$ws.Range("I86").Value = 3933.3333
$ws.Range("J86").Value = 69029
$ws.Range("K86").Value = 3933.3333
$ws.Range("L86").Value = 69029
$ws.Range("M86").Value = -2810.3333
$ws.Range("N86").Value = -71275
$ws.Range("H89").Value = 29971.6
$ws.Range("I89").Value = 3933.3333
$ws.Range("J89").Value = 69029
$ws.Range("K89").Value = 19666.6665
$ws.Range("L89").Value = 345145
$ws.Range("M89").Value = -14050.6665
$ws.Range("N89").Value = -356377
$ws.Range("H105").Value = 843.34784
$ws.Range("I105").Value = 776.3
$ws.Range("K105").Value = 776.3
$ws.Range("M105").Value = 970.7
$ws.Range("H141").Value = 102275.8
$ws.Range("J141").Value = 102275.8
$ws.Range("L141").Value = 102275.8
$ws.Range("N141").Value = -112635.8

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1384.875
$ws.Range("I5").Value = 881.75
$ws.Range("J5").Value = 1888
$ws.Range("K5").Value = 2645.25
$ws.Range("L5").Value = 5664
$ws.Range("M5").Value = -2533.25
$ws.Range("N5").Value = -5888
$ws.Range("H57").Value = 8750
$ws.Range("J57").Value = 8750
$ws.Range("L57").Value = 26250
$ws.Range("N57").Value = -27368
$ws.Range("H131").Value = 782.64514
$ws.Range("I131").Value = 903.3333
$ws.Range("J131").Value = 778.6222
$ws.Range("K131").Value = 2709.9999
$ws.Range("L131").Value = 2335.8666
$ws.Range("M131").Value = 2330.0001
$ws.Range("N131").Value = -12415.8666
$ws.Range("H135").Value = 1384.875
$ws.Range("I135").Value = 881.75
$ws.Range("J135").Value = 1888
$ws.Range("K135").Value = 7935.75
$ws.Range("L135").Value = 16992
$ws.Range("M135").Value = -5400.75
$ws.Range("N135").Value = -22062

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 37.56
$ws.Range("I2").Value = 34.157894
$ws.Range("J2").Value = 48.333332
$ws.Range("K2").Value = 34.157894
$ws.Range("L2").Value = 48.333332
$ws.Range("M2").Value = 78.842106
$ws.Range("N2").Value = -274.333332
$ws.Range("H5").Value = 8666.5
$ws.Range("J5").Value = 13000
$ws.Range("L5").Value = 13000
$ws.Range("N5").Value = -13224

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1042.7435
$ws.Range("I46").Value = 996.4054
$ws.Range("J46").Value = 1900
$ws.Range("K46").Value = 996.4054
$ws.Range("L46").Value = 1900
$ws.Range("M46").Value = -808.4054
$ws.Range("N46").Value = -2276
$ws.Range("H136").Value = 1393.5714
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1559.2142
$ws.Range("I132").Value = 981.1111
$ws.Range("J132").Value = 2599.8
$ws.Range("K132").Value = 2943.3333
$ws.Range("L132").Value = 7799.400000000001
$ws.Range("M132").Value = -413.3332999999998
$ws.Range("N132").Value = -12859.4
